$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.761.59'
$ws.Range('E2').Value = '  -7.55%  '

$ws.Range('D3').Value = '3.688.74'
$ws.Range('E3').Value = '  -7.24%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').Value = '574.05'
$ws.Range('E5').Value = '  -6.93%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '174.40'
$ws.Range('E6').Value = '  +4.67%  '

$ws.Range('D7').Value = '3.679.76'
$ws.Range('E7').Value = '  -7.14%  '

$ws.Range('E8').Value = '  -7.36%  '

$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.03%  '

$ws.Range('D10').Value = '0.712'
$ws.Range('E10').Value = '  -5.78%  '

$ws.Range('D11').Value = '0.165'
$ws.Range('E11').Value = '  -11.72%  '

$ws.Range('D12').Value = '52.75'
$ws.Range('E12').Value = '  -5.09%  '

$ws.Range('D13').Value = '0.0000299'
$ws.Range('E13').Value = '  -11.58%  '

$ws.Range('D14').Value = '10.69'
$ws.Range('E14').Value = '  -3.72%  '

$ws.Range('D15').Value = '4.255.27'

$ws.Range('D16').Value = '3.687.76'
$ws.Range('E16').Value = '  -7.14%  '

$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').Value = '19.37'
$ws.Range('E17').Value = '  -5.34%  '

$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = '0.127'
$ws.Range('E18').Value = '  -3.26%  '

$ws.Range('D19').Value = '1.14'
$ws.Range('E19').Value = '  -8.41%  '

$ws.Range('D20').Value = '12.97'
$ws.Range('E20').Value = '  -7.42%  '

$ws.Range('D21').Value = '67.553.06'
$ws.Range('E21').Value = '  -7.47%  '

$ws.Range('D22').Value = '408.33'
$ws.Range('E22').Value = '  -7.33%  '

$ws.Range('D23').Value = '4.53'
$ws.Range('E23').Value = '  -6.88%  '

$ws.Range('D24').Value = '88.08'
$ws.Range('E24').Value = '  -8.33%  '

$ws.Range('D25').Value = '3.08'
$ws.Range('E25').Value = '  -8.52%  '

$ws.Range('D26').Value = '12.84'
$ws.Range('E26').Value = '  -9.23%  '

$ws.Range('D27').Value = '10.72'
$ws.Range('E27').Value = '  -2.73%  '

$ws.Range('D28').Value = '3.83'
$ws.Range('E28').Value = '  -5.37%  '

$ws.Range('D29').Value = '5.98'
$ws.Range('E29').Value = '  +0.40%  '

$ws.Range('D30').Value = '9.55'
$ws.Range('E30').Value = '  -9.23%  '

$ws.Range('D31').Value = '8.09'
$ws.Range('E31').Value = '  +2.83%  '

$ws.Range('D32').Value = '32.88'
$ws.Range('E32').Value = '  -9.03%  '

$ws.Range('D33').Value = '12.72'
$ws.Range('E33').Value = '  -7.01%  '

$ws.Range('E34').Value = '  -9.16%  '

$ws.Range('D35').Value = '65.35'
$ws.Range('E35').Value = '  -8.02%  '

$ws.Range('D36').Value = '43.89'
$ws.Range('E36').Value = '  -8.00%  '

$ws.Range('D37').Value = '0.0₃0925'
$ws.Range('E37').Value = '  -11.84%  '

$ws.Range('D38').Value = '600.29'
$ws.Range('E38').Value = '  -7.04%  '

$ws.Range('D39').Value = '0.401'
$ws.Range('E39').Value = '  -6.65%  '

$ws.Range('E40').Value = '  +0.11%  '

$ws.Range('E41').Value = '  -0.18%  '

$ws.Range('D42').Value = '3.16'
$ws.Range('E42').Value = '  +6.39%  '

$ws.Range('D43').Value = '0.136'
$ws.Range('E43').Value = '  -6.77%  '

$ws.Range('D44').Value = '3.07'
$ws.Range('E44').Value = '  -10.25%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0440'
$ws.Range('E45').Value = '  -8.70%  '

$ws.Range('D46').Value = '2.62'
$ws.Range('E46').Value = '  +1.97%  '

$ws.Range('D47').Value = '9.44'
$ws.Range('E47').Value = '  -12.10%  '

$ws.Range('D48').Value = '0.135'
$ws.Range('E48').Value = '  -9.12%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.70'
$ws.Range('E49').Value = '  -14.49%  '

$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').Value = '3.12'
$ws.Range('E50').Value = '  -9.23%  '

$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.712.88'
$ws.Range('E51').Value = '  -4.29%  '
